# "Generate Report for Handback"
# Marks the zh-cn and de-de localization rows as handed back: fills in the
# Latest Target File / Latest Handback File / Latest Handback DateTime
# columns, flips the Status from "Ready for handoff" to
# "Handed back: in sync with en-US" (on the Overview rollup + both language
# sheets), and widens the columns whose text just got longer.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"
$sourceMdName = "7f65eb70-c81d-4441-b797-115298bbae1e.md"
$mdHyperlinkAddress = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/080193ffaa8942b5a851167106cb35b596f8168b/e2e/7f65eb70-c81d-4441-b797-115298bbae1e.md"

# Column-width constants:
#  - the "wide" columns (Status on the language sheets, zh-cn/de-de on the
#    Overview sheet) grow to fit the new, longer status text
#  - the handback file-name columns grow to the sheet's standard 40-char
#    file-name column width
$wideWidth = 29.84
$fileColWidth = 39.15

# ----------------------------------------------------------------------
# Overview sheet: both language status columns pick up the new text
# ----------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Columns.Item(5).ColumnWidth = $wideWidth
$overview.Columns.Item(6).ColumnWidth = $wideWidth

# ----------------------------------------------------------------------
# zh-cn sheet: status + handback report for the zh-cn target
# ----------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("J2").Value = $sourceMdName
$zhcn.Range("K2").Value = "7f65eb70-c81d-4441-b797-115298bbae1e.f82a7730ddcc4a03771d16a9eb4bb591ec643bc9.zh-cn.xlf"
$zhcn.Range("L2").Value = "2016-12-13 06:34:13"
$zhcn.Hyperlinks.Add($zhcn.Range("J2"), $mdHyperlinkAddress, "", "", $sourceMdName)
$zhcn.Columns.Item(3).ColumnWidth = $wideWidth
$zhcn.Columns.Item(10).ColumnWidth = $fileColWidth
$zhcn.Columns.Item(11).ColumnWidth = $fileColWidth

# ----------------------------------------------------------------------
# de-de sheet: status + handback report for the de-de target
# ----------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = $newStatus
$dede.Range("J2").Value = $sourceMdName
$dede.Range("K2").Value = "7f65eb70-c81d-4441-b797-115298bbae1e.f82a7730ddcc4a03771d16a9eb4bb591ec643bc9.de-de.xlf"
$dede.Range("L2").Value = "2016-12-13 06:34:31"
$dede.Hyperlinks.Add($dede.Range("J2"), $mdHyperlinkAddress, "", "", $sourceMdName)
$dede.Columns.Item(3).ColumnWidth = $wideWidth
$dede.Columns.Item(10).ColumnWidth = $fileColWidth
$dede.Columns.Item(11).ColumnWidth = $fileColWidth
